$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.801.75'
$ws.Range('E2').Value = '  +1.55%  '
$ws.Range('D3').Value = '2.248.25'
$ws.Range('E3').Value = '  +0.12%  '
$ws.Range('E4').Value = '  +0.58%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.16'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.572'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.01'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.522'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.44'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0806'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.24'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.19%  '
$ws.Range('E13').Value = '  -0.08%  '
$ws.Range('D14').Value = '2.592.20'
$ws.Range('E14').Value = '  +0.20%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.841'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.50%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.244.42'
$ws.Range('E16').Value = '  -0.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.60'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').Value = '44.540.14'
$ws.Range('E18').Value = '  +1.12%  '
$ws.Range('D19').Value = '0.0₃0953'
$ws.Range('E19').Value = '  -1.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.30'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.96'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.96%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.55'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.62'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.96'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.98'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.19%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('E27').Value = '  +3.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.87'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.47'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.96'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '151.10'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.39%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0798'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.97%  '
$ws.Range('E34').Value = '  +0.96%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.08'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.108'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.87%  '
$ws.Range('E37').Value = '  -1.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.86'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.95'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.42'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.81'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0304'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.50%  '
$ws.Range('E43').Value = '  +0.35%  '
$ws.Range('D44').Value = '1.833.67'
$ws.Range('E44').Value = '  +4.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.74'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +14.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.191'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '79.84'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '99.10'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.91'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '69.18'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '54.88'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.91%  '
